$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-20 Thursday" "2025-02-21 Friday"

Replace-Text "674÷8=" "561÷6="
Replace-Text "701÷5=" "174÷7="
Replace-Text "946÷7=" "544÷9="
Replace-Text "194÷2=" "493÷3="
Replace-Text "830÷3=" "562÷6="
Replace-Text "495÷6=" "195÷4="
Replace-Text "437÷5=" "729÷8="
Replace-Text "522÷8=" "783÷4="
Replace-Text "934÷8=" "777÷6="
Replace-Text "509÷8=" "101÷2="
Replace-Text "297÷8=" "942÷9="
Replace-Text "614÷5=" "426÷5="
Replace-Text "234÷4=" "349÷5="
Replace-Text "447÷7=" "861÷4="
Replace-Text "511÷2=" "584÷4="
Replace-Text "600÷6=" "309÷5="
Replace-Text "246÷4=" "621÷2="
Replace-Text "239÷4=" "237÷8="
Replace-Text "128÷2=" "469÷4="
Replace-Text "712÷2=" "465÷2="
Replace-Text "356÷8=" "379÷3="
Replace-Text "848÷7=" "598÷5="
Replace-Text "393÷2=" "285÷7="
Replace-Text "936÷3=" "116÷2="
Replace-Text "507÷4=" "305÷2="
